$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Total" values for B2:B51 (replaces the old "Time" day-fraction values)
$totals = @(7, 7, 8.5, 9.25, 7, 8.5, 7, 8.5, 7, 8.5, 7, 8.5, 7, 7, 7, 7, 7, 7, 8.25, 8.25, 8.25, 8.25, 8.25, 8.25, 9, 9, 9, 9, 9, 9, 7, 7, 7, 7, 7, 7, 7, 8.5, 8.5, 8.5, 8.5, 8.5, 8.5, 7, 7, 7, 7, 7, 7, 7)

# Re-header column B from "Time" to "Total"
$ws.Range("B1").Value = "Total"

# Column B held a time-of-day number format; the new "Total" column is a
# plain number, so clear the inherited formatting before writing data.
$ws.Range("B2:B51").ClearFormats()

for ($i = 0; $i -lt $totals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $totals[$i]
}

# Widen column B so the new header/values are fully visible.
$ws.Columns.Item(2).ColumnWidth = 27.28515625

# Restore the scroll position / selection left by the editing session.
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("B52").Select()
